$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2750.75
$ws.Range("I12").Value = 2900.3333
$ws.Range("J12").Value = 2302
$ws.Range("K12").Value = 2900.3333
$ws.Range("L12").Value = 2302
$ws.Range("M12").Value = -2730.3333
$ws.Range("N12").Value = -2642
$ws.Range("H17").Value = 930.2676
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 914.9857
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 2744.9571
$ws.Range("M17").Value = -5832
$ws.Range("N17").Value = -3080.9571
$ws.Range("H37").Value = 2532.3333
$ws.Range("I37").Value = 999
$ws.Range("J37").Value = 3299
$ws.Range("K37").Value = 2997
$ws.Range("L37").Value = 9897
$ws.Range("M37").Value = -2871
$ws.Range("N37").Value = -10149
$ws.Range("H64").Value = 4053.3845
$ws.Range("I64").Value = 3528.5715
$ws.Range("J64").Value = 4665.6665
$ws.Range("K64").Value = 3528.5715
$ws.Range("L64").Value = 4665.6665
$ws.Range("M64").Value = -3280.5715
$ws.Range("N64").Value = -5161.6665
$ws.Range("H67").Value = 4053.3845
$ws.Range("I67").Value = 3528.5715
$ws.Range("J67").Value = 4665.6665
$ws.Range("K67").Value = 3528.5715
$ws.Range("L67").Value = 4665.6665
$ws.Range("M67").Value = -2670.5715
$ws.Range("N67").Value = -6381.6665
$ws.Range("H74").Value = 4136.4614
$ws.Range("I74").Value = 3815
$ws.Range("J74").Value = 7994
$ws.Range("K74").Value = 3815
$ws.Range("L74").Value = 7994
$ws.Range("M74").Value = -2879
$ws.Range("H77").Value = 4136.4614
$ws.Range("I77").Value = 3815
$ws.Range("J77").Value = 7994
$ws.Range("K77").Value = 19075
$ws.Range("L77").Value = 39970
$ws.Range("M77").Value = -14395
$ws.Range("H106").Value = 173433.83
$ws.Range("I106").Value = 339299.34
$ws.Range("J106").Value = 7568.3335
$ws.Range("K106").Value = 339299.34
$ws.Range("L106").Value = 7568.3335
$ws.Range("M106").Value = -338668.34
$ws.Range("N106").Value = -8830.333500000001
$ws.Range("H116").Value = 6714.077
$ws.Range("I116").Value = 32999
$ws.Range("J116").Value = 4523.6665
$ws.Range("K116").Value = 32999
$ws.Range("L116").Value = 4523.6665
$ws.Range("M116").Value = -29557
$ws.Range("N116").Value = -11407.6665
$ws.Range("H129").Value = 828.9167
$ws.Range("I129").Value = 344.7
$ws.Range("J129").Value = 3250
$ws.Range("K129").Value = 1034.1
$ws.Range("L129").Value = 9750
$ws.Range("M129").Value = 3965.9
$ws.Range("H141").Value = 7379.9585
$ws.Range("I141").Value = 7090.5264
$ws.Range("J141").Value = 8479.799999999999
$ws.Range("K141").Value = 21271.5792
$ws.Range("L141").Value = 25439.4
$ws.Range("M141").Value = -16091.5792
$ws.Range("N141").Value = -35799.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2083.7827
$ws.Range("I2").Value = 2126.3845
$ws.Range("J2").Value = 2028.4
$ws.Range("K2").Value = 2126.3845
$ws.Range("L2").Value = 2028.4
$ws.Range("M2").Value = -2013.3845
$ws.Range("N2").Value = -2254.4
$ws.Range("H12").Value = 3960
$ws.Range("I12").Value = 600
$ws.Range("J12").Value = 9000
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = -427
$ws.Range("H32").Value = 279.13
$ws.Range("I32").Value = 254.08333
$ws.Range("J32").Value = 880.25
$ws.Range("K32").Value = 254.08333
$ws.Range("L32").Value = 880.25
$ws.Range("M32").Value = 32.91667000000001
$ws.Range("H61").Value = 6629.7646
$ws.Range("I61").Value = 6761.5625
$ws.Range("J61").Value = 4521
$ws.Range("K61").Value = 6761.5625
$ws.Range("L61").Value = 4521
$ws.Range("M61").Value = -6549.5625
$ws.Range("N61").Value = -4945
$ws.Range("H74").Value = 9036.5
$ws.Range("I74").Value = 10887.3
$ws.Range("J74").Value = 2867.1667
$ws.Range("K74").Value = 10887.3
$ws.Range("L74").Value = 2867.1667
$ws.Range("M74").Value = -10013.3
$ws.Range("N74").Value = -4615.1667
$ws.Range("H77").Value = 9036.5
$ws.Range("I77").Value = 10887.3
$ws.Range("J77").Value = 2867.1667
$ws.Range("K77").Value = 54436.5
$ws.Range("L77").Value = 14335.8335
$ws.Range("M77").Value = -50068.5
$ws.Range("N77").Value = -23071.8335
$ws.Range("H116").Value = 2083.7827
$ws.Range("I116").Value = 2126.3845
$ws.Range("J116").Value = 2028.4
$ws.Range("K116").Value = 2126.3845
$ws.Range("L116").Value = 2028.4
$ws.Range("M116").Value = 167.6154999999999
$ws.Range("N116").Value = -6616.4
$ws.Range("H132").Value = 3682.9412
$ws.Range("I132").Value = 3157.68
$ws.Range("J132").Value = 5142
$ws.Range("K132").Value = 9473.039999999999
$ws.Range("L132").Value = 15426
$ws.Range("M132").Value = -6943.039999999999
$ws.Range("N132").Value = -20486
$ws.Range("H136").Value = 6629.7646
$ws.Range("I136").Value = 6761.5625
$ws.Range("J136").Value = 4521
$ws.Range("K136").Value = 20284.6875
$ws.Range("L136").Value = 13563
$ws.Range("M136").Value = -17734.6875
$ws.Range("N136").Value = -18663

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 149712
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 149712
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 149712
$ws.Range("N2").Value = -149938
$ws.Range("H3").Value = 2083.7827
$ws.Range("I3").Value = 2126.3845
$ws.Range("J3").Value = 2028.4
$ws.Range("K3").Value = 2126.3845
$ws.Range("L3").Value = 2028.4
$ws.Range("M3").Value = -2012.3845
$ws.Range("N3").Value = -2256.4
$ws.Range("H20").Value = 2113.6
$ws.Range("I20").Value = 2058
$ws.Range("J20").Value = 2231.75
$ws.Range("K20").Value = 2058
$ws.Range("L20").Value = 2231.75
$ws.Range("M20").Value = -1811
$ws.Range("N20").Value = -2725.75
$ws.Range("H81").Value = 154149
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 154149
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 154149
$ws.Range("N81").Value = -156271
$ws.Range("H84").Value = 154149
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 154149
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 462447
$ws.Range("N84").Value = -473055
$ws.Range("H94").Value = 2177.3635
$ws.Range("I94").Value = 2100.9
$ws.Range("J94").Value = 2942
$ws.Range("K94").Value = 2100.9
$ws.Range("L94").Value = 2942
$ws.Range("M94").Value = -1649.9
$ws.Range("N94").Value = -3844
$ws.Range("H134").Value = 5450.577
$ws.Range("I134").Value = 5414.1304
$ws.Range("J134").Value = 5730
$ws.Range("K134").Value = 16242.3912
$ws.Range("L134").Value = 17190
$ws.Range("M134").Value = -13707.3912
$ws.Range("N134").Value = -22260

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1679.8889
$ws.Range("I16").Value = 717
$ws.Range("J16").Value = 2050.2307
$ws.Range("K16").Value = 717
$ws.Range("L16").Value = 2050.2307
$ws.Range("M16").Value = -430
$ws.Range("H41").Value = 11388.214
$ws.Range("I41").Value = 8043.5
$ws.Range("J41").Value = 19750
$ws.Range("K41").Value = 8043.5
$ws.Range("L41").Value = 19750
$ws.Range("M41").Value = -7615.5
$ws.Range("N41").Value = -20606
$ws.Range("H50").Value = 50000
$ws.Range("I50").Value = 42500
$ws.Range("J50").Value = 80000
$ws.Range("K50").Value = 42500
$ws.Range("L50").Value = 80000
$ws.Range("M50").Value = -41875
$ws.Range("N50").Value = -81250
$ws.Range("H58").Value = 5903.5293
$ws.Range("I58").Value = 5788.973
$ws.Range("J58").Value = 6206.2856
$ws.Range("K58").Value = 5788.973
$ws.Range("L58").Value = 6206.2856
$ws.Range("M58").Value = -5585.973
$ws.Range("N58").Value = -6612.2856
$ws.Range("H74").Value = 79545.45
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 79545.45
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 79545.45
$ws.Range("N74").Value = -81293.45
$ws.Range("H77").Value = 79545.45
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 79545.45
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 238636.35
$ws.Range("N77").Value = -247372.35
$ws.Range("H105").Value = 3969.9167
$ws.Range("I105").Value = 2222
$ws.Range("J105").Value = 4552.5557
$ws.Range("K105").Value = 2222
$ws.Range("L105").Value = 4552.5557
$ws.Range("M105").Value = -475
$ws.Range("N105").Value = -8046.5557
$ws.Range("H107").Value = 5277
$ws.Range("I107").Value = 4999
$ws.Range("J107").Value = 5555
$ws.Range("K107").Value = 4999
$ws.Range("L107").Value = 5555
$ws.Range("M107").Value = -3079
$ws.Range("N107").Value = -9395
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = $null
$ws.Range("H113").Value = 1679.8889
$ws.Range("I113").Value = 717
$ws.Range("J113").Value = 2050.2307
$ws.Range("K113").Value = 717
$ws.Range("L113").Value = 2050.2307
$ws.Range("M113").Value = 1453
$ws.Range("H131").Value = 75394.8
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 75394.8
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 75394.8
$ws.Range("M131").Value = $null
$ws.Range("N131").Value = -85474.8
$ws.Range("H132").Value = 7800.3335
$ws.Range("I132").Value = 3305.3333
$ws.Range("J132").Value = 21285.334
$ws.Range("K132").Value = 9915.999899999999
$ws.Range("L132").Value = 63856.00199999999
$ws.Range("M132").Value = -7385.999899999999
$ws.Range("H134").Value = 7055.033
$ws.Range("I134").Value = 7591.0454
$ws.Range("J134").Value = 5581
$ws.Range("K134").Value = 22773.1362
$ws.Range("L134").Value = 16743
$ws.Range("M134").Value = -20238.1362
$ws.Range("N134").Value = -21813
$ws.Range("H136").Value = 5903.5293
$ws.Range("I136").Value = 5788.973
$ws.Range("J136").Value = 6206.2856
$ws.Range("K136").Value = 17366.919
$ws.Range("L136").Value = 18618.8568
$ws.Range("M136").Value = -14816.919
$ws.Range("N136").Value = -23718.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 64.15385000000001
$ws.Range("I2").Value = 23.5
$ws.Range("J2").Value = 89.5625
$ws.Range("K2").Value = 141
$ws.Range("L2").Value = 537.375
$ws.Range("M2").Value = -28
$ws.Range("N2").Value = -763.375
$ws.Range("H7").Value = 224.25
$ws.Range("I7").Value = 170
$ws.Range("J7").Value = 314.66666
$ws.Range("K7").Value = 510
$ws.Range("L7").Value = 943.9999799999999
$ws.Range("M7").Value = -398
$ws.Range("N7").Value = -1167.99998
$ws.Range("H12").Value = 162.07143
$ws.Range("I12").Value = 121
$ws.Range("J12").Value = 178.5
$ws.Range("K12").Value = 363
$ws.Range("L12").Value = 535.5
$ws.Range("M12").Value = -190
$ws.Range("H50").Value = 2521.5386
$ws.Range("I50").Value = 359.2
$ws.Range("J50").Value = 3873
$ws.Range("K50").Value = 1077.6
$ws.Range("L50").Value = 11619
$ws.Range("M50").Value = -596.5999999999999
$ws.Range("N50").Value = -12581
$ws.Range("H53").Value = 2521.5386
$ws.Range("I53").Value = 359.2
$ws.Range("J53").Value = 3873
$ws.Range("K53").Value = 1077.6
$ws.Range("L53").Value = 11619
$ws.Range("M53").Value = -596.5999999999999
$ws.Range("N53").Value = -12581
$ws.Range("H59").Value = 10
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 10
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 30
$ws.Range("M59").Value = $null
$ws.Range("N59").Value = -1110
$ws.Range("H74").Value = 15000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 45000
$ws.Range("M74").Value = $null
$ws.Range("N74").Value = -47122
$ws.Range("H77").Value = 15000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 135000
$ws.Range("M77").Value = $null
$ws.Range("N77").Value = -145608
$ws.Range("H87").Value = 3000
$ws.Range("I87").Value = 3000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 9000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -7752
$ws.Range("H90").Value = 3000
$ws.Range("I90").Value = 3000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 27000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -20760
$ws.Range("I94").Value = 9986.5
$ws.Range("J94").Value = 139999920
$ws.Range("K94").Value = 29959.5
$ws.Range("L94").Value = 419999760
$ws.Range("M94").Value = -29283.5
$ws.Range("N94").Value = -420001112
$ws.Range("H115").Value = 1526.0952
$ws.Range("I115").Value = 277.0909
$ws.Range("J115").Value = 2900
$ws.Range("K115").Value = 831.2727
$ws.Range("L115").Value = 8700
$ws.Range("M115").Value = 343.7273
$ws.Range("N115").Value = -11050
$ws.Range("H116").Value = 4251722
$ws.Range("I116").Value = 4637787.5
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 13913362.5
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = -13909920.5
$ws.Range("H118").Value = 3350
$ws.Range("I118").Value = 3350
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 10050
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -8807
$ws.Range("N118").Value = $null
$ws.Range("H119").Value = 4964.5
$ws.Range("I119").Value = 4964.5
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 14893.5
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -10055.5
$ws.Range("N119").Value = $null
$ws.Range("H120").Value = 1197.5
$ws.Range("I120").Value = 1197.5
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 3592.5
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = 1245.5
$ws.Range("H132").Value = 5046.4165
$ws.Range("I132").Value = 3980.5715
$ws.Range("J132").Value = 6538.6
$ws.Range("K132").Value = 35825.1435
$ws.Range("L132").Value = 58847.4
$ws.Range("M132").Value = -33295.1435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20912518
$ws.Range("I80").Value = 30669934
$ws.Range("J80").Value = 3767.2856
$ws.Range("K80").Value = 30669934
$ws.Range("L80").Value = 3767.2856
$ws.Range("M80").Value = -30668936
$ws.Range("N80").Value = -5763.2856
$ws.Range("H83").Value = 20912518
$ws.Range("I83").Value = 30669934
$ws.Range("J83").Value = 3767.2856
$ws.Range("K83").Value = 153349670
$ws.Range("L83").Value = 18836.428
$ws.Range("M83").Value = -153344678
$ws.Range("N83").Value = -28820.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 857.1111
$ws.Range("I22").Value = 866.2
$ws.Range("J22").Value = 845.75
$ws.Range("K22").Value = 866.2
$ws.Range("L22").Value = 845.75
$ws.Range("M22").Value = -571.2
$ws.Range("N22").Value = -1435.75
$ws.Range("H27").Value = 857.1111
$ws.Range("I27").Value = 866.2
$ws.Range("J27").Value = 845.75
$ws.Range("K27").Value = 866.2
$ws.Range("L27").Value = 845.75
$ws.Range("M27").Value = -759.2
$ws.Range("N27").Value = -1059.75
$ws.Range("H35").Value = 11095
$ws.Range("I35").Value = 3038
$ws.Range("J35").Value = 16466.334
$ws.Range("K35").Value = 3038
$ws.Range("L35").Value = 16466.334
$ws.Range("M35").Value = -2702
$ws.Range("N35").Value = -17138.334
$ws.Range("H40").Value = 4440.643
$ws.Range("I40").Value = 4455.6895
$ws.Range("J40").Value = 4407.077
$ws.Range("K40").Value = 4455.6895
$ws.Range("L40").Value = 4407.077
$ws.Range("M40").Value = -4319.6895
$ws.Range("H68").Value = 3584.6155
$ws.Range("I68").Value = 1356.5217
$ws.Range("J68").Value = 20666.666
$ws.Range("K68").Value = 1356.5217
$ws.Range("L68").Value = 20666.666
$ws.Range("M68").Value = -607.5217
$ws.Range("H71").Value = 3584.6155
$ws.Range("I71").Value = 1356.5217
$ws.Range("J71").Value = 20666.666
$ws.Range("K71").Value = 6782.6085
$ws.Range("L71").Value = 103333.33
$ws.Range("M71").Value = -3038.6085
$ws.Range("H132").Value = 37671.844
$ws.Range("I132").Value = 47858.793
$ws.Range("J132").Value = 7111
$ws.Range("K132").Value = 143576.379
$ws.Range("L132").Value = 21333
$ws.Range("M132").Value = -141046.379

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = $null
$ws.Range("H62").Value = 11766.333
$ws.Range("I62").Value = 10700
$ws.Range("J62").Value = 12832.667
$ws.Range("K62").Value = 10700
$ws.Range("L62").Value = 12832.667
$ws.Range("M62").Value = -10076
$ws.Range("H65").Value = 11766.333
$ws.Range("I65").Value = 10700
$ws.Range("J65").Value = 12832.667
$ws.Range("K65").Value = 53500
$ws.Range("L65").Value = 64163.335
$ws.Range("M65").Value = -50380
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null
$ws.Range("H122").Value = 14733.482
$ws.Range("I122").Value = 14135.115
$ws.Range("J122").Value = 19919.334
$ws.Range("K122").Value = 42405.345
$ws.Range("L122").Value = 59758.00199999999
$ws.Range("M122").Value = -39955.345
$ws.Range("H126").Value = 9119.944
$ws.Range("I126").Value = 5110.727
$ws.Range("J126").Value = 15420.143
$ws.Range("K126").Value = 15332.181
$ws.Range("L126").Value = 46260.429
$ws.Range("M126").Value = -12862.181
$ws.Range("N126").Value = -51200.429
